# Updates crypto market data cells to the latest scraped snapshot.
# D-column "Price" cells hold text that often LOOKS numeric (e.g. "68.451.86"
# using '.' as a thousands separator, or "561.87"); Excel's COM Value setter
# auto-coerces numeric-looking strings to real numbers, which would silently
# change both the stored type and the effective value. To keep those cells
# as plain text (matching the source data), we briefly force the cell to
# Text format ("@") before assigning the value, then clear the format again
# so the cell's style reverts to the sheet's default (matching the
# unstyled inline-string cells produced by the original export).
function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "68.451.86"
$ws.Range("E2").Value = "  -1.69%  "
Set-TextValue $ws "D3" "2.456.48"
$ws.Range("E3").Value = "  -1.84%  "
Set-TextValue $ws "D5" "561.87"
$ws.Range("E5").Value = "  -2.49%  "
Set-TextValue $ws "D6" "163.87"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.60%  "
Set-TextValue $ws "D9" "2.456.37"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("E12").Value = "  -4.76%  "
Set-TextValue $ws "D13" "4.83"
$ws.Range("E13").Value = "  -2.14%  "
Set-TextValue $ws "D14" "2.917.99"
$ws.Range("E14").Value = "  -1.39%  "
Set-TextValue $ws "D15" "68.406.83"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("E16").Value = "  -3.57%  "
Set-TextValue $ws "D17" "23.48"
$ws.Range("E17").Value = "  -5.07%  "
Set-TextValue $ws "D18" "2.551.68"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("E19").Value = "  -2.07%  "
Set-TextValue $ws "D20" "343.31"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("E21").Value = "  -4.10%  "
Set-TextValue $ws "D22" "3.81"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -3.27%  "
Set-TextValue $ws "D25" "67.97"
$ws.Range("E25").Value = "  -3.51%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D26" "3.73"
$ws.Range("E26").Value = "  -5.68%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws "D27" "1.06"
$ws.Range("E27").Value = "  +6.04%  "
$ws.Range("E28").Value = "  -1.29%  "
Set-TextValue $ws "D29" "8.21"
$ws.Range("E29").Value = "  -6.24%  "
Set-TextValue $ws "D30" "0.0₃0838"
$ws.Range("E30").Value = "  -6.04%  "
Set-TextValue $ws "D31" "7.29"
$ws.Range("E31").Value = "  -6.72%  "
Set-TextValue $ws "D32" "3.29"
$ws.Range("E32").Value = "  +120.39%  "
$ws.Range("E33").Value = "  -2.68%  "
Set-TextValue $ws "D34" "434.28"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("E35").Value = "  +0.01%  "
Set-TextValue $ws "D36" "1.68"
$ws.Range("E36").Value = "  -2.88%  "
Set-TextValue $ws "D37" "157.58"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D39" "0.110"
$ws.Range("E39").Value = "  -5.47%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D40" "1.00"
$ws.Range("E40").Value = "  +0.00%  "
Set-TextValue $ws "D41" "17.91"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("E42").Value = "  -3.11%  "
Set-TextValue $ws "D43" "4.48"
$ws.Range("E43").Value = "  -4.22%  "
Set-TextValue $ws "D44" "1.52"
$ws.Range("E44").Value = "  -4.53%  "
Set-TextValue $ws "D45" "1.10"
$ws.Range("E45").Value = "  +1.19%  "
Set-TextValue $ws "D46" "2.07"
$ws.Range("E46").Value = "  -5.76%  "
Set-TextValue $ws "D47" "134.42"
$ws.Range("E47").Value = "  -5.01%  "
$ws.Range("E48").Value = "  -3.11%  "
Set-TextValue $ws "D49" "0.0717"
$ws.Range("E49").Value = "  -2.12%  "
Set-TextValue $ws "D50" "0.486"
$ws.Range("E50").Value = "  -6.41%  "
Set-TextValue $ws "D51" "0.562"
$ws.Range("E51").Value = "  -2.79%  "
